$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("AE2").Value = 29
$ws.Range("AH2").Value = 34
$ws.Range("AT2").Value = 2.25

# Row 3
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5

# Row 4
$ws.Range("BD4").Value = 126

# Row 6
$ws.Range("G6").Value = 1.55
$ws.Range("I6").Value = 6
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.33
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 2.1
$ws.Range("U6").Value = 1.73
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 8
$ws.Range("AC6").Value = 13
$ws.Range("AE6").Value = 15
$ws.Range("AJ6").Value = 51
$ws.Range("AM6").Value = 201
$ws.Range("AP6").Value = 17
$ws.Range("AS6").Value = 101
$ws.Range("AU6").Value = 8
$ws.Range("BA6").Value = 101

$wb.Save()
